# Workbook-level edit script
# Reverts a set of test-data tweaks: rotate numeric suffixes (31 -> 33),
# fix up some accented-character strings to their ASCII-safe equivalents,
# adjust a couple of hyperlinks, and update selected cells / active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Shared "identity" values live on Sheet4 (A2/C2/E2). Every other sheet
#    (Sheet2, Sheet3, Sheet5, Sheet6, Sheet7, Sheet8) pulls them in via
#    formulas, so changing them here cascades everywhere automatically.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("A2").Value = "tavalinetont33"
$ws4.Range("C2").Value = "puhtaloom33"
$ws4.Range("E2").Value = "filmweird33"

# ---------------------------------------------------------------------
# 2. Direct text fixes (these are plain values, not formulas).
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Sheet5")
$ws5.Range("G2").Value = "Jalgpalliass"
$ws5.Range("I2").Value = "selentest20@mailinator.com"
$ws5.Range("I3").Value = "testimeauto20@mailinator.com"
$ws5.Range("I4").Value = "vahekonto20@mailinator.com"

$ws6 = $wb.Worksheets.Item("Sheet6")
$ws6.Range("D2").Value = "Kas see on tavaline kysitlus"
$ws6.Range("E2").Value = "Tundub taiesti tavaline"
$ws6.Range("J2").Value = "6nnenumbrid on mul"
$ws6.Range("M2").Value = "Yks"
$ws6.Range("D4").Value = "Kas liiga palju jarge on halb?"
$ws6.Range("E4").Value = "Ei, kuid ju siis toovad raha sisse"

$ws8 = $wb.Worksheets.Item("Sheet8")
$ws8.Range("F2").Value = "Väravaid oskab lyya, kuid kaitsa ei oska"

# ---------------------------------------------------------------------
# 3. Hyperlinks on Sheet5: drop the I3 / I4 mail links, keep J2/J3/J4
#    (all three pointing at selentest@hotmail.com), renumbering the
#    relationship ids in the process.
# ---------------------------------------------------------------------
$ws5.Range("I3").Hyperlinks.Delete()
$ws5.Hyperlinks.Add($ws5.Range("J2"), "mailto:selentest@hotmail.com", "", "", "selentest@hotmail.com")
$ws5.Hyperlinks.Add($ws5.Range("J3"), "mailto:selentest@hotmail.com", "", "", "selentest@hotmail.com")
$ws5.Hyperlinks.Add($ws5.Range("J4"), "mailto:selentest@hotmail.com", "", "", "selentest@hotmail.com")

# ---------------------------------------------------------------------
# 4. Selection / active-tab bookkeeping.
# ---------------------------------------------------------------------
$ws4.Range("D38").Select()

$ws6.Range("I28").Select()

$ws8.Range("G2").Select()

# Activating Sheet5 makes it the active tab and marks tabSelected on it,
# clearing the previous tabSelected flag on Sheet9 automatically.
$ws5.Activate()
$ws5.Range("I4").Select()

# Sheet9's selection itself doesn't move (still I4); only tabSelected
# changes, which activating Sheet5 above already took care of.
